$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.160.42"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "1.603.49"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.89"
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("E8").Value = "  +3.13%  "
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +4.86%  "
$ws.Range("D12").Value = "1.825.95"
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("D13").Value = "1.605.09"
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.512"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "26.136.97"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.55"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "205.32"
$ws.Range("E20").Value = "  +11.34%  "
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.83"
$ws.Range("E24").Value = "  +10.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.87"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("E28").Value = "  +3.27%  "
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0164"
$ws.Range("E36").Value = "  +10.67%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.115.59"
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").Value = "1.738.74"
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.95"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  +5.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.45"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.409"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0922"
$ws.Range("E51").Value = "  -13.17%  "
